# Weekly crime data refresh: shift report window forward one week
# (Volume 31 Number 48, week of 11/25/2024-12/1/2024 -> Number 49, week of 12/2/2024-12/8/2024)
# and update the Week to Date / 28 Day / Year to Date / 2 Year crime stat table (rows 14-33).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/week number and report date range) ---
$ws.Range("A8").Value = "Volume 31   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/2/2024  Through  12/8/2024"

# --- Crime Complaints table (rows 14-33), columns C:N ---
$cellValues = @{
    "C14" = "'0"
    "D14" = "'0"
    "E14" = "'***.*"
    "F14" = "'0"
    "G14" = 1
    "H14" = -100
    "I14" = 3
    "J14" = 5
    "K14" = -40
    "L14" = -40
    "M14" = "'***.*"
    "N14" = -81.25
    "C15" = 1
    "D15" = "'0"
    "E15" = "'***.*"
    "F15" = 1
    "G15" = 1
    "H15" = 0
    "I15" = 17
    "J15" = 17
    "K15" = 0
    "L15" = 30.76923076923
    "M15" = -10.526315789473
    "N15" = -15
    "C16" = "'0"
    "D16" = 3
    "E16" = -100
    "F16" = 12
    "G16" = 14
    "H16" = -14.285714285714
    "I16" = 202
    "J16" = 192
    "K16" = 5.208333333333
    "L16" = -3.34928229665
    "M16" = -9.821428571428
    "N16" = -53.348729792147
    "C17" = 1
    "D17" = 3
    "E17" = -66.666666666666
    "F17" = 15
    "G17" = 18
    "H17" = -16.666666666666
    "I17" = 315
    "J17" = 275
    "K17" = 14.545454545454
    "L17" = 27.016129032258
    "M17" = 83.13953488372
    "N17" = 42.533936651583
    "C18" = "'0"
    "D18" = "'0"
    "E18" = "'***.*"
    "F18" = 6
    "G18" = 2
    "H18" = 200
    "I18" = 168
    "J18" = 123
    "K18" = 36.585365853658
    "L18" = 68
    "M18" = -41.868512110726
    "N18" = -78.125
    "C19" = 10
    "D19" = 12
    "E19" = -16.666666666666
    "F19" = 44
    "G19" = 54
    "H19" = -18.518518518518
    "I19" = 655
    "J19" = 632
    "K19" = 3.639240506329
    "L19" = 17.383512544802
    "M19" = 43.640350877193
    "N19" = 64.987405541561
    "C20" = 6
    "D20" = 9
    "E20" = -33.333333333333
    "F20" = 35
    "G20" = 24
    "H20" = 45.833333333333
    "I20" = 409
    "J20" = 478
    "K20" = -14.435146443514
    "L20" = 43.006993006993
    "M20" = 105.527638190955
    "N20" = -80.270139893873
    "C21" = 18
    "D21" = 27
    "E21" = -33.333333333333
    "F21" = 113
    "G21" = 114
    "H21" = -0.877192982456
    "I21" = 1769
    "J21" = 1722
    "K21" = 2.729384436701
    "L21" = 24.665257223396
    "M21" = 30.169242089771
    "N21" = -54.964358452138
    "C22" = "'0"
    "D22" = 1
    "E22" = -100
    "F22" = 1
    "G22" = 3
    "H22" = -66.666666666666
    "I22" = 8
    "J22" = 21
    "K22" = -61.904761904761
    "L22" = -46.666666666666
    "M22" = -38.461538461538
    "N22" = "'***.*"
    "C23" = "'0"
    "D23" = 1
    "E23" = -100
    "F23" = 1
    "G23" = 1
    "H23" = 0
    "I23" = 54
    "J23" = 55
    "K23" = -1.818181818181
    "L23" = 74.193548387096
    "M23" = 5.882352941176
    "N23" = "'***.*"
    "C24" = 12
    "D24" = 27
    "E24" = -55.555555555555
    "F24" = 132
    "G24" = 121
    "H24" = 9.090909090909
    "I24" = 1329
    "J24" = 1326
    "K24" = 0.226244343891
    "L24" = 19.407008086253
    "M24" = -1.918819188191
    "N24" = "'***.*"
    "C25" = 6
    "D25" = 9
    "E25" = -33.333333333333
    "F25" = 56
    "G25" = 78
    "H25" = -28.205128205128
    "I25" = 708
    "J25" = 771
    "K25" = -8.17120622568
    "L25" = 33.333333333333
    "M25" = "'***.*"
    "N25" = "'***.*"
    "C26" = 6
    "D26" = 3
    "E26" = 100
    "F26" = 31
    "G26" = 30
    "H26" = 3.333333333333
    "I26" = 515
    "J26" = 472
    "K26" = 9.110169491525
    "L26" = 20.327102803738
    "M26" = 25
    "N26" = "'***.*"
    "C27" = 1
    "D27" = "'0"
    "E27" = "'***.*"
    "F27" = 1
    "G27" = 3
    "H27" = -66.666666666666
    "I27" = 25
    "J27" = 32
    "K27" = -21.875
    "L27" = -7.407407407407
    "M27" = "'***.*"
    "N27" = "'***.*"
    "C28" = 1
    "D28" = "'0"
    "E28" = "'***.*"
    "F28" = 3
    "G28" = 2
    "H28" = 50
    "I28" = 61
    "J28" = 41
    "K28" = 48.780487804878
    "L28" = 29.787234042553
    "M28" = "'***.*"
    "N28" = "'***.*"
    "C29" = "'0"
    "D29" = "'0"
    "E29" = "'***.*"
    "F29" = "'0"
    "G29" = 1
    "H29" = -100
    "I29" = 4
    "J29" = 10
    "K29" = -60
    "L29" = -69.230769230769
    "M29" = -42.857142857142
    "N29" = -81.818181818181
    "C30" = "'0"
    "D30" = "'0"
    "E30" = "'***.*"
    "F30" = "'0"
    "G30" = 1
    "H30" = -100
    "I30" = 4
    "J30" = 8
    "K30" = -50
    "L30" = -55.555555555555
    "M30" = -42.857142857142
    "N30" = -81.818181818181
    "C31" = "'0"
    "D31" = "'0"
    "E31" = "'***.*"
    "F31" = "'0"
    "G31" = "'0"
    "H31" = "'***.*"
    "I31" = 1
    "J31" = 4
    "K31" = -75
    "L31" = -83.333333333333
    "M31" = "'***.*"
    "N31" = "'***.*"
    "C33" = 1
    "D33" = "'0"
    "E33" = "'***.*"
    "F33" = 1
    "G33" = "'0"
    "H33" = "'***.*"
    "I33" = 8
    "J33" = 4
    "K33" = 100
    "L33" = 14.285714285714
    "M33" = "'***.*"
    "N33" = "'***.*"
}

foreach ($ref in $cellValues.Keys) {
    $ws.Range($ref).Value = $cellValues[$ref]
}
